{"js": "// \"Updated What the Function rules to say Python 3\"\n//\n// The rules document refers to the Python programming language twice;\n// both mentions should become \"Python 3\". We locate every exact,\n// case-sensitive occurrence of the word \"Python\" in the body and insert\n// \" 3\" immediately after it. This reproduces the only real content\n// change in the target revision (the rest of the underlying diff is\n// just Word re-chunking unrelated runs on save, with no visible text\n// change).\n\nconst body = context.document.body;\n\nconst pythonMatches = body.search(\"Python\", {\n  matchCase: true,\n  matchWholeWord: true\n});\npythonMatches.load(\"text\");\nawait context.sync();\n\nfor (let i = 0; i < pythonMatches.items.length; i++) {\n  pythonMatches.items[i].insertText(\" 3\", Word.InsertLocation.after);\n}\nawait context.sync();\n\n// Word also keeps a \"_GoBack\" bookmark marking the location of the most\n// recent edit (so Shift+F5 can return the cursor there). In the target\n// revision it moved from the very first paragraph to sit inside the\n// \"basic Python 3 syntax\" sentence, right after \"incl\" (where the\n// author's cursor ended up while typing). Re-create that so the\n// bookmark position matches the authored document as closely as\n// possible.\ntry {\n  context.document.deleteBookmark(\"_GoBack\");\n  await context.sync();\n} catch (e) {\n  // no-op if it doesn't exist / host doesn't support the call\n}\n\nconst anchor = body.search(\"syntax, incl\", { matchCase: true });\nawait context.sync();\n\nif (anchor.items.length > 0) {\n  const afterAnchor = anchor.items[0].getRange(Word.RangeLocation.after);\n  afterAnchor.insertBookmark(\"_GoBack\");\n  await context.sync();\n}\n", "ps1": "# \"Updated What the Function rules to say Python 3\"\n#\n# The rules document refers to the Python programming language twice;\n# both mentions should become \"Python 3\". Find & Replace every exact,\n# whole-word, case-sensitive occurrence of \"Python\" with \"Python 3\" in\n# the document body. That's the only real content change in the target\n# revision (the remainder of the underlying OOXML diff is just Word\n# re-chunking unrelated runs on save, with no visible text change).\n\n$d = $word.ActiveDocument\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Execute(\"Python\", $true, $true, $false, $false, $false, $true, 1, $false, \"Python 3\", 2)\n\n# Word also keeps a \"_GoBack\" bookmark marking the location of the most\n# recent edit (so Shift+F5 can return the cursor there). In the target\n# revision it moved from the very first paragraph to sit inside the\n# \"basic Python 3 syntax\" sentence, right after \"incl\" (where the\n# author's cursor ended up while typing). Re-create that so the\n# bookmark position matches the authored document as closely as\n# possible.\nif ($d.Bookmarks.Exists(\"_GoBack\")) {\n    $d.Bookmarks(\"_GoBack\").Delete()\n}\n\n$find2 = $d.Content.Find\n$find2.ClearFormatting()\n$find2.Text = \"syntax, incl\"\n$find2.Forward = $true\n$find2.Wrap = 1\n$find2.MatchCase = $true\n$found = $find2.Execute()\n\nif ($found) {\n    $anchorRange = $find2.Parent\n    $target = $d.Range($anchorRange.End, $anchorRange.End)\n    $d.Bookmarks.Add(\"_GoBack\", $target)\n}\n"}
